$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.509852886199951
$ws.Range("B1").Value = 3.571891069412231
$ws.Range("C1").Value = 3.157248258590698
$ws.Range("D1").Value = 3.959398984909058
$ws.Range("E1").Value = 5.434846878051758
